# "Generate Report for Archive"
# The status "Ready for handoff" moves to "In Translation" everywhere it is
# used (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2 -- all four cells share
# the same shared string), and the Status/locale columns that held that text
# are re-sized (narrower, since the new text is shorter than the old text).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
